$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-RowValues {
    param($ws, $row, $values)
    $col = 2
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

# Update existing rows 218-223 with corrected figures
Set-RowValues $ws 218 @(7177,3647,3283,2130,1400,5350,1653,660,2235,623,1243)
Set-RowValues $ws 219 @(7294,4510,4254,1918,867,5735,1836,709,2230,848,1291)
Set-RowValues $ws 220 @(7761,4953,4554,2170,638,6882,2314,921,2610,997,1523)
Set-RowValues $ws 221 @(8160,5067,4598,2513,580,5986,1906,849,2377,814,1380)
Set-RowValues $ws 222 @(7621,5029,4663,2204,388,6693,2166,933,2631,914,1540)
Set-RowValues $ws 223 @(7453,4803,4375,2217,432,6381,2035,774,2603,940,1366)

# Add new row 224 for 01-07-2021 (force text so it isn't parsed as a date,
# then clear the format back to the default so no extra cell style lingers)
$cellA224 = $ws.Cells.Item(224, 1)
$cellA224.NumberFormat = "@"
$cellA224.Value = "01-07-2021"
$cellA224.ClearFormats()
Set-RowValues $ws 224 @(7943,5144,4584,2379,419,7339,2493,1009,2979,1092,1491)
